# Generate Report for handoff
# Update row 3 (b.md) status/handoff info across Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: b.md row (row 3) status for each locale ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) ---
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-01-26 09:34:15"

$zhLink = $wsZhCn.Hyperlinks.Item(1)
for ($i = 1; $i -le $wsZhCn.Hyperlinks.Count; $i++) {
    $lnk = $wsZhCn.Hyperlinks.Item($i)
    if ($lnk.Range.Address -eq "$3" -or $lnk.Range.Row -eq 3) {
    }
}

# --- de-de sheet: b.md row (row 3) ---
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-01-26 09:34:27"
